$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6 (ReasonForCallValue): ElementType changes from "Drop Down" to "Text Box",
# and the XPath locator in column G is removed.
$ws.Range("F6").Value = "Text Box"
$ws.Range("G6").Value = ""

# Row 7 (ScreeningNarrative): ElementType changes from "Text Area" to "Text Box",
# and the XPath locator in column G is removed.
$ws.Range("F7").Value = "Text Box"
$ws.Range("G7").Value = ""

# Row 8 (CallerType): XPath locator in column G is removed.
$ws.Range("G8").Value = ""

# Row 9 (CallBackRequiredDrpDwn): ElementType changes from "Text Box" to "Drop Down".
$ws.Range("F9").Value = "Drop Down"

# Row 10 (CallBackRequiredValue): XPath locator in column G is removed.
$ws.Range("G10").Value = ""

# Row 11 (SaveAndProceed): XPath locator in column G is removed.
$ws.Range("G11").Value = ""

# Rows 6-8 no longer need the taller, wrapped row height now that the long
# XPath text is gone - revert them back to the sheet's default row height.
$null = $ws.Rows("6:8").AutoFit()

# The active selection moved to F9 in the saved workbook.
$null = $ws.Range("F9").Select()
